$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The mall list had three duplicate/unwanted hair-dye rows (Pink Hair Dye,
# Light Green Hair dye, White Hair dye) sitting between "Super Hair Gel"
# and "Light Blue hair dye". Remove those three rows entirely - the rows
# below shift up, the shared strings they used are no longer referenced,
# and the shared formulas in column E recompute against the new C values.
$ws.Rows("54:56").Delete()

# Scroll the view down to where the list now ends and leave the same
# relative selection the author left it on.
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 45
$ws.Range("F50").Select()
